# Commit: "Diagramas de robustez y registro de tiempo al hacerlos"
# Register progress + time spent on the robustness diagrams for
# use cases CU-14, CU-15, CU-16 and CU-21 (rows 59, 61, 63, 65 of
# the "Casos de Uso" sheet) and update the sheet views.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# The status column (F) moves from "Por iniciar" to "En proceso" for the
# four tasks that now have work logged against them.
$ws.Range("F59").Value = "En proceso"
$ws.Range("F61").Value = "En proceso"
$ws.Range("F63").Value = "En proceso"
$ws.Range("F65").Value = "En proceso"

# Hours consumed registered for Día 5 (column T) on CU-14 / CU-15.
$ws.Range("T59").Value = 0.5
$ws.Range("T61").Value = 0.5

# Hours consumed registered for Día 6 (column W) on CU-16 / CU-21.
$ws.Range("W63").Value = 1
$ws.Range("W65").Value = 1

# Update the frozen-pane view state to where the user left off editing.
$csUso = $excel.Windows.Item(1)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 19
$ws.Range("F65").Select()

$wsInstr = $wb.Worksheets.Item("Instructivo")
$wsInstr.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsInstr.Range("C9").Select()

$ws.Activate()
